$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$shp = $s.Shapes.Item(1)
$tbl = $shp.Table
Write-Host "Rows.Count before:" $tbl.Rows.Count
while ($tbl.Rows.Count -gt 2) {
    $tbl.Rows.Item($tbl.Rows.Count).Delete()
}
Write-Host "Rows.Count after:" $tbl.Rows.Count
Write-Host "Shape.Height read1:" $shp.Height
Write-Host "Shape.Height read2:" $shp.Height
Write-Host "Shape.Height read3:" $shp.Height
$shp.Width = 500
Write-Host "Shape.Height after width set:" $shp.Height
